# Updated RAD New Tax Return Amount Due Test Data: refresh the recorded
# execution timestamps in the "Date" column (B) of the NewTaxReturn sheet
# to reflect the latest test run (Tue Mar 05 2024 01:10:21 .. 01:25:47 EST).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewTaxReturn")

$dates = @(
    "Tue Mar 05 01:10:21 EST 2024",
    "Tue Mar 05 01:10:38 EST 2024",
    "Tue Mar 05 01:10:56 EST 2024",
    "Tue Mar 05 01:11:14 EST 2024",
    "Tue Mar 05 01:11:31 EST 2024",
    "Tue Mar 05 01:11:49 EST 2024",
    "Tue Mar 05 01:12:06 EST 2024",
    "Tue Mar 05 01:12:24 EST 2024",
    "Tue Mar 05 01:12:41 EST 2024",
    "Tue Mar 05 01:12:59 EST 2024",
    "Tue Mar 05 01:13:17 EST 2024",
    "Tue Mar 05 01:13:37 EST 2024",
    "Tue Mar 05 01:13:54 EST 2024",
    "Tue Mar 05 01:14:12 EST 2024",
    "Tue Mar 05 01:14:30 EST 2024",
    "Tue Mar 05 01:14:47 EST 2024",
    "Tue Mar 05 01:15:07 EST 2024",
    "Tue Mar 05 01:15:26 EST 2024",
    "Tue Mar 05 01:15:43 EST 2024",
    "Tue Mar 05 01:16:11 EST 2024",
    "Tue Mar 05 01:16:30 EST 2024",
    "Tue Mar 05 01:16:49 EST 2024",
    "Tue Mar 05 01:17:06 EST 2024",
    "Tue Mar 05 01:17:25 EST 2024",
    "Tue Mar 05 01:17:44 EST 2024",
    "Tue Mar 05 01:18:04 EST 2024",
    "Tue Mar 05 01:18:21 EST 2024",
    "Tue Mar 05 01:18:41 EST 2024",
    "Tue Mar 05 01:18:59 EST 2024",
    "Tue Mar 05 01:19:19 EST 2024",
    "Tue Mar 05 01:19:36 EST 2024",
    "Tue Mar 05 01:19:54 EST 2024",
    "Tue Mar 05 01:20:13 EST 2024",
    "Tue Mar 05 01:20:32 EST 2024",
    "Tue Mar 05 01:20:49 EST 2024",
    "Tue Mar 05 01:21:09 EST 2024",
    "Tue Mar 05 01:21:28 EST 2024",
    "Tue Mar 05 01:21:47 EST 2024",
    "Tue Mar 05 01:22:05 EST 2024",
    "Tue Mar 05 01:22:24 EST 2024",
    "Tue Mar 05 01:22:43 EST 2024",
    "Tue Mar 05 01:23:02 EST 2024",
    "Tue Mar 05 01:23:20 EST 2024",
    "Tue Mar 05 01:23:37 EST 2024",
    "Tue Mar 05 01:23:56 EST 2024",
    "Tue Mar 05 01:24:15 EST 2024",
    "Tue Mar 05 01:24:32 EST 2024",
    "Tue Mar 05 01:24:51 EST 2024",
    "Tue Mar 05 01:25:10 EST 2024",
    "Tue Mar 05 01:25:29 EST 2024",
    "Tue Mar 05 01:25:47 EST 2024"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}
